$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.764.82'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '2.675.97'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.61'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '2.675.34'
$ws.Range('E9').Value = '  -1.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.144'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.159'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.92'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').Value = '3.164.03'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000185'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.08%  '
$ws.Range('D17').Value = '67.670.48'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').Value = '2.675.38'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.77'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.76'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '364.22'
$ws.Range('D21').ClearFormats()
$ws.Range('E22').Value = '  -3.22%  '
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('E24').Value = '  -3.94%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.03'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.16'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.49%  '
$ws.Range('D28').Value = '2.815.95'
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('E29').Value = '  -2.69%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '556.09'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -6.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.04'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('E33').Value = '  -3.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.93'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.96%  '
$ws.Range('E35').Value = '  -1.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  -4.67%  '
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '155.49'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('E40').Value = '  -1.91%  '
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('E42').Value = '  -4.31%  '
$ws.Range('E44').Value = '  -6.39%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.30'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('E47').Value = '  -5.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.591'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '153.63'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.87'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('E51').Value = '  -3.42%  '
